# Remove the two "private Expression parse...Expr()" paragraphs from the
# "Relevant Parser Methods" slide (slide 9), leaving the surrounding
# paragraphs (and the trailing blank paragraphs) untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)   # "Content Placeholder 2"
$tr = $shp.TextFrame.TextRange

# Paragraphs (1-based):
#   1: private InitialDecl parseStringTypeDecl()
#   2: private Variable parseVariable()
#   3: private Expression parseIndexExpr()      <- remove
#   4: private Expression parseFieldExpr()      <- remove
#   5: (blank)
#   6: (blank)
$firstPara = $tr.Paragraphs(3, 1)
$lastPara  = $tr.Paragraphs(4, 1)

$startPos = $firstPara.Start
$endPos   = $lastPara.Start + $lastPara.Length
$len      = $endPos - $startPos

$tr.Characters($startPos, $len).Delete()
